{"js": "const replacements = [\n  [\"2025-08-28 Thursday\", \"2025-08-29 Friday\"],\n  [\"749\u00f76=124, 5\", \"333\u00f79=37, 0\"],\n  [\"115\u00f79=12, 7\", \"379\u00f78=47, 3\"],\n  [\"673\u00f74=168, 1\", \"997\u00f73=332, 1\"],\n  [\"847\u00f79=94, 1\", \"967\u00f75=193, 2\"],\n  [\"702\u00f76=117, 0\", \"726\u00f75=145, 1\"],\n  [\"176\u00f78=22, 0\", \"850\u00f75=170, 0\"],\n  [\"342\u00f76=57, 0\", \"647\u00f75=129, 2\"],\n  [\"766\u00f77=109, 3\", \"445\u00f73=148, 1\"],\n  [\"668\u00f73=222, 2\", \"751\u00f72=375, 1\"],\n  [\"270\u00f72=135, 0\", \"941\u00f78=117, 5\"],\n  [\"973\u00f77=139, 0\", \"455\u00f72=227, 1\"],\n  [\"926\u00f78=115, 6\", \"493\u00f78=61, 5\"],\n  [\"975\u00f74=243, 3\", \"562\u00f79=62, 4\"],\n  [\"728\u00f75=145, 3\", \"379\u00f79=42, 1\"],\n  [\"786\u00f75=157, 1\", \"866\u00f72=433, 0\"],\n  [\"952\u00f79=105, 7\", \"453\u00f73=151, 0\"],\n  [\"817\u00f75=163, 2\", \"520\u00f74=130, 0\"],\n  [\"229\u00f76=38, 1\", \"986\u00f74=246, 2\"],\n  [\"802\u00f76=133, 4\", \"876\u00f75=175, 1\"],\n  [\"466\u00f73=155, 1\", \"177\u00f79=19, 6\"],\n  [\"839\u00f72=419, 1\", \"595\u00f73=198, 1\"],\n  [\"935\u00f74=233, 3\", \"228\u00f74=57, 0\"],\n  [\"729\u00f78=91, 1\", \"405\u00f76=67, 3\"],\n  [\"788\u00f74=197, 0\", \"753\u00f79=83, 6\"],\n  [\"565\u00f72=282, 1\", \"924\u00f79=102, 6\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2025-08-28 Thursday', '2025-08-29 Friday')\n    ,@('749\u00f76=124, 5', '333\u00f79=37, 0')\n    ,@('115\u00f79=12, 7', '379\u00f78=47, 3')\n    ,@('673\u00f74=168, 1', '997\u00f73=332, 1')\n    ,@('847\u00f79=94, 1', '967\u00f75=193, 2')\n    ,@('702\u00f76=117, 0', '726\u00f75=145, 1')\n    ,@('176\u00f78=22, 0', '850\u00f75=170, 0')\n    ,@('342\u00f76=57, 0', '647\u00f75=129, 2')\n    ,@('766\u00f77=109, 3', '445\u00f73=148, 1')\n    ,@('668\u00f73=222, 2', '751\u00f72=375, 1')\n    ,@('270\u00f72=135, 0', '941\u00f78=117, 5')\n    ,@('973\u00f77=139, 0', '455\u00f72=227, 1')\n    ,@('926\u00f78=115, 6', '493\u00f78=61, 5')\n    ,@('975\u00f74=243, 3', '562\u00f79=62, 4')\n    ,@('728\u00f75=145, 3', '379\u00f79=42, 1')\n    ,@('786\u00f75=157, 1', '866\u00f72=433, 0')\n    ,@('952\u00f79=105, 7', '453\u00f73=151, 0')\n    ,@('817\u00f75=163, 2', '520\u00f74=130, 0')\n    ,@('229\u00f76=38, 1', '986\u00f74=246, 2')\n    ,@('802\u00f76=133, 4', '876\u00f75=175, 1')\n    ,@('466\u00f73=155, 1', '177\u00f79=19, 6')\n    ,@('839\u00f72=419, 1', '595\u00f73=198, 1')\n    ,@('935\u00f74=233, 3', '228\u00f74=57, 0')\n    ,@('729\u00f78=91, 1', '405\u00f76=67, 3')\n    ,@('788\u00f74=197, 0', '753\u00f79=83, 6')\n    ,@('565\u00f72=282, 1', '924\u00f79=102, 6')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
